$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data (text-formatted cells)
$cell = $ws.Cells.Item(2, 4)
$cell.NumberFormat = "@"
$cell.Value = "68.195.64"
$cell = $ws.Cells.Item(2, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -7.00%  "

$cell = $ws.Cells.Item(3, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.709.60"
$cell = $ws.Cells.Item(3, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -6.91%  "

$cell = $ws.Cells.Item(4, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.06%  "

$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = "584.47"
$cell = $ws.Cells.Item(5, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -5.19%  "

$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = "176.74"
$cell = $ws.Cells.Item(6, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +5.71%  "

$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.700.24"
$cell = $ws.Cells.Item(7, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -7.00%  "

$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.635"
$cell = $ws.Cells.Item(8, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -7.08%  "

$cell = $ws.Cells.Item(9, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.25%  "

$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.716"
$cell = $ws.Cells.Item(10, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -5.06%  "

$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.166"
$cell = $ws.Cells.Item(11, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -10.80%  "

$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = "52.94"
$cell = $ws.Cells.Item(12, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -6.44%  "

$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0000302"
$cell = $ws.Cells.Item(13, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -10.82%  "

$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = "10.67"
$cell = $ws.Cells.Item(14, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -4.06%  "

$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = "4.305.94"
$cell = $ws.Cells.Item(15, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -6.86%  "

$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.734.89"
$cell = $ws.Cells.Item(16, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -6.54%  "

$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = "19.42"
$cell = $ws.Cells.Item(17, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -5.41%  "

$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.127"
$cell = $ws.Cells.Item(18, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -3.13%  "

$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = "13.05"
$cell = $ws.Cells.Item(19, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -7.71%  "

$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.14"
$cell = $ws.Cells.Item(20, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -8.30%  "

$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = "68.098.06"
$cell = $ws.Cells.Item(21, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -7.03%  "

$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = "410.80"
$cell = $ws.Cells.Item(22, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -6.80%  "

$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = "4.66"
$cell = $ws.Cells.Item(23, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -4.62%  "

$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = "89.06"
$cell = $ws.Cells.Item(24, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -7.01%  "

$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.09"
$cell = $ws.Cells.Item(25, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -8.37%  "

$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = "12.88"
$cell = $ws.Cells.Item(26, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -9.16%  "

$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = "10.77"
$cell = $ws.Cells.Item(27, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -3.55%  "

$cell = $ws.Cells.Item(28, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -5.79%  "

$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = "5.95"
$cell = $ws.Cells.Item(29, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.09%  "

$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = "9.61"
$cell = $ws.Cells.Item(30, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -8.50%  "

$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = "8.07"
$cell = $ws.Cells.Item(31, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +3.30%  "

$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = "33.07"
$cell = $ws.Cells.Item(32, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -8.45%  "

$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = "12.80"
$cell = $ws.Cells.Item(33, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -6.57%  "

$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = "44.82"
$cell = $ws.Cells.Item(34, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -5.92%  "

$cell = $ws.Cells.Item(35, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -8.79%  "

$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = "66.17"
$cell = $ws.Cells.Item(36, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -7.17%  "

$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0935"
$cell = $ws.Cells.Item(37, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -7.66%  "

$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = "599.00"
$cell = $ws.Cells.Item(38, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -6.23%  "

$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.406"
$cell = $ws.Cells.Item(39, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -6.27%  "

$cell = $ws.Cells.Item(40, 2)
$cell.NumberFormat = "@"
$cell.Value = "dogwifhat"
$cell = $ws.Cells.Item(40, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.38"
$cell = $ws.Cells.Item(40, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +16.09%  "

$cell = $ws.Cells.Item(41, 2)
$cell.NumberFormat = "@"
$cell.Value = "Dai"
$cell = $ws.Cells.Item(41, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell = $ws.Cells.Item(41, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.01%  "

$cell = $ws.Cells.Item(42, 2)
$cell.NumberFormat = "@"
$cell.Value = "FirstDigitalUSD"
$cell = $ws.Cells.Item(42, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell = $ws.Cells.Item(42, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.12%  "

$cell = $ws.Cells.Item(43, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -6.71%  "

$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.08"
$cell = $ws.Cells.Item(44, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -10.38%  "

$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0442"
$cell = $ws.Cells.Item(45, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -8.76%  "

$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = "9.59"
$cell = $ws.Cells.Item(46, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -13.03%  "

$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.60"
$cell = $ws.Cells.Item(47, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.51%  "

$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.135"
$cell = $ws.Cells.Item(48, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -9.19%  "

$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.754.73"
$cell = $ws.Cells.Item(49, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -3.44%  "

$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.15"
$cell = $ws.Cells.Item(50, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -8.15%  "

$cell = $ws.Cells.Item(51, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -4.87%  "
